# Apply updated cryptocurrency price/volume figures (and the VeChain/Maker
# row swap) to sheet1, matching the Sun May  5 09:54:38 UTC 2024 data refresh.
#
# Price-column (D) updates are written with a temporary Text number format
# so Excel keeps them as literal strings (e.g. "587.43") instead of
# auto-converting them to floating point numbers; the format is reset back
# to the default "Normal" style right after the value is set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.777.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.143.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.144.67'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  +6.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.663.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.143.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.574.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  +1.79%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.65%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.03'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0855'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("E38").Value = '  -4.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '440.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.78'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.86%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.921.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0371'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.277'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.77'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.39%  '

Write-Output 'Applied cryptos update.'
